$d = $word.ActiveDocument

# --- 1) "The basics dataset:" bullet -----------------------------------
# Merge the trailing " " / "Contains the basic information..." /
# " including language, type, attributes of the movie" runs into a
# single run by re-writing the (already contiguous) text in place.
$rng1 = $d.Content
$rng1.Find.Execute(
    " Contains the basic information about each movie including language, type, attributes of the movie",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " Contains the basic information about each movie including language, type, attributes of the movie",
    2) | Out-Null

# --- 2) "The titles dataset:" bullet ------------------------------------
# Same treatment: collapse the leading space run + the detail run into one.
$rng2 = $d.Content
$rng2.Find.Execute(
    " Contains extra details including title, release year, runtime, genre, adult-rating",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " Contains extra details including title, release year, runtime, genre, adult-rating",
    2) | Out-Null

# --- 3) New "Data Source" paragraph -------------------------------------
# Inserted after the "To focus on a particular question..." paragraph and
# before the "Disclaimer:" paragraph.
$focusPara = $d.Paragraphs.Item(8)
$focusPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Item(9)
$newRng = $newPara.Range
# Exclude the paragraph mark itself so subsequent inserts stay inside
# this paragraph instead of spilling into the next one.
$newRng.End = $newRng.End - 1

$newRng.InsertAfter("Datasets are very large in size and therefore are not uploaded to ") | Out-Null
$newRng.Collapse(0)
$newRng.InsertAfter("Github") | Out-Null
$newRng.Collapse(0)
$newRng.InsertAfter(". The datasets can be found at: ") | Out-Null
$newRng.Collapse(0)
$newRng.InsertAfter("https://www.imdb.com/interfaces/") | Out-Null
